# Fill in gas-estimation data rows 15-19 (previously blank B/C inputs, so
# D/E/F evaluated to "") and extend the table with new rows 20-35
# (formulas only, inputs stay blank) so the sheet now spans A1:F35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Give the new rows (18-35) the same cell formatting as row 14, the
#     last fully-populated row, before writing any values/formulas into
#     them. Rows 15-17 already existed with the right formatting. ---
$ws.Range("A14:F14").Copy()
$ws.Range("A18:F35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 15-19: fill in the B (tx cost) / C (gas price) inputs ---
$data = @{
    15 = @(3077295, 41)
    16 = @(2944022, 38)
    17 = @(2944022, 38)
    18 = @(2691648, 38)
    19 = @(2607152, 38)
}

foreach ($r in 15..19) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
}
# Rows 15-17 already carry the D/E/F formulas (shared with the D5:D17 /
# E6:E17 / F4:F17 groups) from the original sheet, so simply populating
# B/C above is enough for them to recompute.

# --- D/E/F formulas, new rows 18-35: one Range.Formula assignment per
#     column so Excel groups them into fresh shared-formula blocks. ---
$ws.Range("D18:D35").Formula = '=IF(AND(B18<>"",C18<>""),B18*C18,"")'
$ws.Range("E18:E35").Formula = '=IF(D18<>"",IF(C18=C17,D18-D17,"--"),"")'
$ws.Range("F18:F35").Formula = '=IF(D18<>"",D18*0.000000001,"")'

# --- View state: scroll so row 4 is the top row, selection on H8 ---
$ws.Range("H8").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
